# MakeVersionInfo: a couple minor changes
#
# Update the "Date Placeholder" text shown on every slide, every slide
# layout, and the slide master from 2024-07-24 to 2024-08-17.

$p = $ppt.ActivePresentation

$oldDate = "2024-07-24"
$newDate = "2024-08-17"
$ppPlaceholderDate = 16

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)

        $isDatePlaceholder = $false
        if ($sh.Type -eq 14) {
            try {
                if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDatePlaceholder = $true
                }
            } catch {
                $isDatePlaceholder = $false
            }
        }

        if ($isDatePlaceholder -and $sh.HasTextFrame) {
            $tf = $sh.TextFrame
            if ($tf.HasText -and ($tf.TextRange.Text -eq $oldDate)) {
                $tf.TextRange.Text = $newDate
            }
        }
    }
}

# 1. Every slide.
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    Update-DatePlaceholders $s.Shapes
}

# 2. Every slide layout (CustomLayouts, reached through the Design).
$design = $p.Designs.Item(1)
$layouts = $design.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $lay = $layouts.Item($li)
    Update-DatePlaceholders $lay.Shapes
}

# 3. The slide master itself.
Update-DatePlaceholders $p.SlideMaster.Shapes
